# Product Backlog: mark stories #46 and #47 as DONE by cyan-highlighting
# their "User Story" cell text (matching the existing convention used for
# other completed rows in the table, e.g. story #48's green/cyan rows).
#
# Note: this runtime's `Range.HighlightColorIndex` setter does not apply to
# the target Range, so highlighting is applied via Find/Replace: searching
# for the literal text and replacing it with itself while
# `Find.Replacement.Highlight` is set to the WdColorIndex for cyan (3).
# This reliably writes <w:highlight w:val="cyan"/> into the run(s) that
# contain the matched text.
$cyan = 3

function Set-CyanHighlight([string]$text) {
    $find = $word.ActiveDocument.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $text
    $find.Replacement.Text = $text
    $find.Replacement.Highlight = $cyan
    return $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
}

$d = $word.ActiveDocument

# --- Story #46: "edit user accounts" -------------------------------------
# Single run, text itself is unchanged -- just add the cyan highlight.
$story46 = "As an admin, I should be able to edit user accounts so that I can update user accounts in case of login issues"
Set-CyanHighlight $story46 | Out-Null

# --- Story #47: "delete/remove user accounts" -----------------------------
# This paragraph's text is split across two runs with a <w:bookmarkStart/
# w:bookmarkEnd w:name="_GoBack"/> sitting between them. In the target
# revision the run split point shifts earlier (from right before "nt
# System." to right before "er accounts...") while the bookmark stays put
# structurally between the two runs, and both runs get the cyan highlight.
#
# A straight Find/Replace across the whole paragraph would delete the
# bookmark (it falls inside the replaced span), so instead the text is
# first re-split at the new boundary via direct Range edits (which do not
# disturb the bookmark, since those ranges stay within a single run), and
# only then is the cyan highlight applied to each of the two resulting
# runs separately.

$fullText = "As an admin, I should be able to delete/remove user accounts so that I can remove user accounts that are no longer have access in the PUP Medical Clinic Record Management System."
$firstRunText = "As an admin, I should be able to delete/remove user accounts so that I can remove user accounts that are no longer have access in the PUP Medical Clinic Record Manageme"

# Locate the paragraph (Find.Execute matches straight through the bookmark).
$paraRange = $d.Content
$paraRange.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraStart = $paraRange.Start
$paraEnd = $paraRange.End

# Locate the current run1/run2 boundary (i.e. where the bookmark currently sits).
$run1Range = $d.Content
$run1Range.Find.Execute($firstRunText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$run1End = $run1Range.Start + $firstRunText.Length

$run2Text = $d.Range($run1End, $paraEnd).Text

# New split point, matching the target revision exactly.
$newSplitOffset = 84
$newSplitAbs = $paraStart + $newSplitOffset
$movedText = $firstRunText.Substring($newSplitOffset)

# Rewrite run2 first (this Range sits entirely inside the old run2, so it
# cannot disturb the bookmark which sits just before it).
$run2Range = $d.Range($run1End, $paraEnd)
$run2Range.Text = $movedText + $run2Text

# Now trim run1's tail down to the new split point (this Range sits
# entirely inside run1 too, so the bookmark -- still located at $run1End,
# now beyond the trimmed text -- is left untouched).
$tailOfRun1 = $d.Range($newSplitAbs, $run1End)
$tailOfRun1.Text = ""

# Finally, highlight the two (now correctly split) runs.
$newRun1Text = "As an admin, I should be able to delete/remove user accounts so that I can remove us"
$newRun2Text = "er accounts that are no longer have access in the PUP Medical Clinic Record Management System."
Set-CyanHighlight $newRun1Text | Out-Null
Set-CyanHighlight $newRun2Text | Out-Null
